$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.334207
$ws.Range("H2").Value = 1.002621
$ws.Range("I2").Value = 0.07226389998643547
$ws.Range("J2").Value = 0.07226389998643548
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.002279333333333333
$ws.Range("N2").Value = 0.006838
$ws.Range("Q2").Value = 0.0007617691553333332
$ws.Range("R2").Value = 0.006855922398000001
$ws.Range("S2").Value = 0.07226389998643547
$ws.Range("T2").Value = 0.07226389998643548
$ws.Range("I3").Value = 0.5366421328200824
$ws.Range("J3").Value = 0.5366421328200826
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.002279333333333333
$ws.Range("N3").Value = 0.006838
$ws.Range("Q3").Value = 0.005657007500444445
$ws.Range("R3").Value = 0.050913067504
$ws.Range("S3").Value = 0.5366421328200824
$ws.Range("T3").Value = 0.5366421328200826
$ws.Range("I4").Value = 0.3910939671934819
$ws.Range("J4").Value = 0.391093967193482
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.002279333333333333
$ws.Range("N4").Value = 0.006838
$ws.Range("Q4").Value = 0.004122713015777778
$ws.Range("R4").Value = 0.037104417142
$ws.Range("S4").Value = 0.3910939671934819
$ws.Range("T4").Value = 0.391093967193482

Write-Output "Applied TPM update to Podxl2-Sell sheet"
